$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.09215
$ws.Range("N2").Value = 0.1843
$ws.Range("O2").Value = 0.01023932159503481
$ws.Range("P2").Value = 0.01017879497979401
$ws.Range("Q2").Value = 0.0004079787666666667
$ws.Range("R2").Value = 0.0024478726
$ws.Range("S2").Value = 0.01023932159503481
$ws.Range("T2").Value = 0.01017879497979401

# Update row 3 values
$ws.Range("O3").Value = 0.01189268776136058
$ws.Range("P3").Value = 0.01773358166721151
$ws.Range("Q3").Value = 0.000473856010888889
$ws.Range("R3").Value = 0.004264704098000001
$ws.Range("S3").Value = 0.01189268776136058
$ws.Range("T3").Value = 0.01773358166721151

# Update row 4 values
$ws.Range("M4").Value = 8.80044
$ws.Range("N4").Value = 17.60088
$ws.Range("O4").Value = 0.9778679906436047
$ws.Range("P4").Value = 0.9720876233529945
$ws.Range("Q4").Value = 0.03896248136
$ws.Range("R4").Value = 0.23377488816
$ws.Range("S4").Value = 0.9778679906436047
$ws.Range("T4").Value = 0.9720876233529945

# Delete rows 5 and 6
$ws.Range("A5:T6").EntireRow.Delete()
